$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Kreiranje porudžbine" paragraph - the old _GoBack bookmark split the
# run between "...prikazane su i slik" and "e, kako bi...". Re-merge that text
# into a single run (the Find/Replace spans the bookmark gap, which removes
# the bookmark and normalizes the two adjacent identically-formatted runs
# into one).
# ---------------------------------------------------------------------------
$rng1 = $d.Content
[void]$rng1.Find.Execute("prikazane su i slike, kako bi konobari", $true, $false, $false, $false, $false, $true, 1, $false, "prikazane su i slike, kako bi konobari", 2)

# ---------------------------------------------------------------------------
# Edit 2: "Izmena porudžbine" paragraph - color the sentence "Sistem pri
# izmeni porudžbine obaveštava kuvare i šankere o novim podacima" red,
# splitting it out of the surrounding run (the trailing ". " stays black).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
[void]$rng2.Find.Execute("Sistem pri izmeni porudžbine obaveštava kuvare i šankere o novim podacima", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Font.Color = 255

# ---------------------------------------------------------------------------
# Edit 3: move the _GoBack bookmark to the end of the (empty, single-space)
# paragraph right before "Nakon završetka pripreme jela...". A perfectly
# zero-width range placed exactly at a paragraph-mark offset resolves to the
# wrong side of the boundary, so we nudge it: temporarily extend the run by
# one throw-away character, drop the bookmark immediately before that
# character (a safe, non-boundary offset), then remove the character again.
# Adding a bookmark named "_GoBack" also automatically relocates any
# existing bookmark of that name, so the old one (removed already by Edit 1)
# doesn't need separate cleanup.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
[void]$rng3.Find.Execute("Nakon završetka pripreme jela, on šalje notifikaciju da je jelo spremno.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $rng3.Paragraphs(1)
$prevPara = $targetPara.Previous()
$prevRange = $prevPara.Range
[void]$prevRange.MoveEnd(1, -1)
$prevRange.Collapse(0)
$insertPos = $prevRange.Start
$prevRange.InsertAfter("Z")
$bmTarget = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmTarget)
$tempCharRange = $d.Range($insertPos, $insertPos + 1)
$tempCharRange.Delete()
